$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates ---
# (set in the same order the strings first appear in the saved workbook so
#  the shared-string table order matches the target file)
$ws.Range("L4").Value = "edit Use Case"
$ws.Range("L5").Value = "specify entities and relation"
$ws.Range("L8").Value = "refine use case- sequence and activity diagram"
$ws.Range("L7").Value = "constrain diagram- reduce all contexts, like scenario carrying person, putsize of human"
$ws.Range("L6").Value = "specify stakeholder- humans that be rescued, operation"
$ws.Range("E5").Value = "3) Use Case Diagram"

# --- Column width changes ---
$ws.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12.833333333333334
$ws.Columns.Item(5).ColumnWidth = 34.833333333333336
$ws.Columns.Item(6).ColumnWidth = 29.666666666666668
$ws.Columns.Item(7).ColumnWidth = 16.666666666666668
$ws.Columns.Item(9).ColumnWidth = 13.666666666666666
$ws.Columns.Item(10).ColumnWidth = 21.666666666666668
$ws.Columns.Item(11).ColumnWidth = 17.666666666666668
$ws.Columns.Item(12).ColumnWidth = 70.0
for ($c = 13; $c -le 19; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 10.666666666666666
}

# --- View / scroll / selection changes ---
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E5").Select()

# --- Application window size/position (best-effort; mirrors author's resize) ---
$excel.ActiveWindow.Left = -108
$excel.ActiveWindow.Top = -108
$excel.ActiveWindow.Width = 1008
$excel.ActiveWindow.Height = 546
